# Update the cached "datetimeFigureOut" footer-date field from 08/06/2020
# to 12/06/2020 on the slide master and on every slide layout, then rename
# the "specimenRequirements" action label to "specimenRequested" on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "08/06/2020"
$newDate = "12/06/2020"
$ppPlaceholderDate = 16

# --- Slide master: date placeholder -----------------------------------
$master = $p.SlideMaster
for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $sh = $master.Shapes.Item($si)
    if ($sh.HasTextFrame) {
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Every slide layout: date placeholder ------------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# --- Slide 1: rename the "specimenRequirements" run --------------------
# Note: a run that happens to be the last run in its paragraph comes back
# from .Runs()/.Paragraphs() with a trailing CR (paragraph-mark) character,
# so trim before comparing and only assign the clean replacement text back.
$slide = $p.Slides.Item(1)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $sh = $slide.Shapes.Item($si)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            $runCount = $para.Runs().Count
            for ($ri = 1; $ri -le $runCount; $ri++) {
                $run = $para.Runs($ri, 1)
                $runText = $run.Text.TrimEnd("`r", "`n")
                if ($runText -eq "specimenRequirements ") {
                    $run.Text = "specimenRequested "
                }
            }
        }
    }
}
